$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 12 ("products") to make room for
# the new "logged_in" field that documents the accounts table's new
# boolean column. All rows from 12 downward shift down by one.
$ws.Rows.Item(11).EntireRow.Insert()

# Populate the newly inserted row 11 with the new field description.
$ws.Range("A11").Value = "logged_in"
$ws.Range("B11").Value = "Will be a boolean, 0 or 1"

# Restore the selection used while editing the sheet.
$ws.Range("B13").Select() | Out-Null
